$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1605.1111
$ws.Range("I12").Value = 1064.2858
$ws.Range("K12").Value = 1064.2858
$ws.Range("M12").Value = -894.2858000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1683.2142
$ws.Range("I28").Value = 1641.8462
$ws.Range("J28").Value = 2221
$ws.Range("K28").Value = 1641.8462
$ws.Range("L28").Value = 2221
$ws.Range("M28").Value = -1156.8462
$ws.Range("N28").Value = -3191

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 773.5
$ws.Range("I32").Value = 773.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 773.5
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -447.5
$ws.Range("N32").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3404.7273
$ws.Range("I43").Value = 2244
$ws.Range("K43").Value = 2244
$ws.Range("M43").Value = -2175

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7951.7856
$ws.Range("I62").Value = 7613.846
$ws.Range("K62").Value = 7613.846
$ws.Range("M62").Value = -6989.846

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 7951.7856
$ws.Range("I65").Value = 7613.846
$ws.Range("K65").Value = 38069.23
$ws.Range("M65").Value = -34949.23

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1924.2667
$ws.Range("I98").Value = 1943.909
$ws.Range("K98").Value = 1943.909
$ws.Range("M98").Value = -445.9090000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 51493.75
$ws.Range("J105").Value = 51493.75
$ws.Range("L105").Value = 51493.75
$ws.Range("N105").Value = -58481.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 2645.5454
$ws.Range("I107").Value = 2840.2
$ws.Range("K107").Value = 2840.2
$ws.Range("M107").Value = -920.1999999999998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2390.5
$ws.Range("I112").Value = 2949.75
$ws.Range("K112").Value = 8849.25
$ws.Range("M112").Value = -7741.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1924.2667
$ws.Range("I122").Value = 1943.909
$ws.Range("K122").Value = 5831.727000000001
$ws.Range("M122").Value = -3381.727000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1350
$ws.Range("I129").Value = 1350
$ws.Range("K129").Value = 4050
$ws.Range("M129").Value = 950

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 5549
$ws.Range("I141").Value = 5080.4736
$ws.Range("K141").Value = 15241.4208
$ws.Range("M141").Value = -10061.4208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 700
$ws.Range("I25").Value = 700
$ws.Range("K25").Value = 700
$ws.Range("M25").Value = -298

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23810736
$ws.Range("I32").Value = 25000948
$ws.Range("K32").Value = 25000948
$ws.Range("M32").Value = -25000661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2157.9666
$ws.Range("I74").Value = 2201.36
$ws.Range("J74").Value = 1941
$ws.Range("K74").Value = 2201.36
$ws.Range("L74").Value = 1941
$ws.Range("M74").Value = -1327.36
$ws.Range("N74").Value = -3689

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2157.9666
$ws.Range("I77").Value = 2201.36
$ws.Range("J77").Value = 1941
$ws.Range("K77").Value = 11006.8
$ws.Range("L77").Value = 9705
$ws.Range("M77").Value = -6638.800000000001
$ws.Range("N77").Value = -18441

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 94819.5
$ws.Range("J121").Value = 94819.5
$ws.Range("L121").Value = 94819.5
$ws.Range("N121").Value = -98313.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = -127

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 95990
$ws.Range("J122").Value = 95990
$ws.Range("L122").Value = 95990
$ws.Range("N122").Value = -105790

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5522.0938
$ws.Range("I31").Value = 4001.8
$ws.Range("K31").Value = 4001.8
$ws.Range("M31").Value = -3706.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5522.0938
$ws.Range("I34").Value = 4001.8
$ws.Range("K34").Value = 4001.8
$ws.Range("M34").Value = -3799.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 70612
$ws.Range("J87").Value = 70612
$ws.Range("L87").Value = 70612
$ws.Range("N87").Value = -72984

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H90").Value = 70612
$ws.Range("J90").Value = 70612
$ws.Range("L90").Value = 211836
$ws.Range("N90").Value = -223692

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 52427
$ws.Range("I103").Value = 17250
$ws.Range("J103").Value = 99329.664
$ws.Range("K103").Value = 17250
$ws.Range("L103").Value = 99329.664
$ws.Range("M103").Value = -16078
$ws.Range("N103").Value = -101673.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 3364.3333
$ws.Range("J32").Value = 4999
$ws.Range("L32").Value = 14997
$ws.Range("N32").Value = -15563

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1734.7142
$ws.Range("I131").Value = 1532.5
$ws.Range("K131").Value = 4597.5
$ws.Range("M131").Value = 442.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 4452.1875
$ws.Range("I133").Value = 4203.1816
$ws.Range("K133").Value = 12609.5448
$ws.Range("M133").Value = -7549.5448

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 56833.668
$ws.Range("J119").Value = 56833.668
$ws.Range("L119").Value = 56833.668
$ws.Range("N119").Value = -66509.66800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2667.1667
$ws.Range("I122").Value = 2631.3333
$ws.Range("J122").Value = 2774.6667
$ws.Range("K122").Value = 7893.999899999999
$ws.Range("L122").Value = 8324.000100000001
$ws.Range("M122").Value = -5443.999899999999
$ws.Range("N122").Value = -13224.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3333.12
$ws.Range("I132").Value = 3284.5
$ws.Range("K132").Value = 9853.5
$ws.Range("M132").Value = -7323.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3641.75
$ws.Range("I46").Value = 3000.75
$ws.Range("J46").Value = 3748.5833
$ws.Range("K46").Value = 3000.75
$ws.Range("L46").Value = 3748.5833
$ws.Range("M46").Value = -2812.75
$ws.Range("N46").Value = -4124.5833

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 7000
$ws.Range("I64").Value = 5000
$ws.Range("K64").Value = 5000
$ws.Range("M64").Value = -4775

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H67").Value = 7000
$ws.Range("I67").Value = 5000
$ws.Range("K67").Value = 5000
$ws.Range("M67").Value = -4220

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2239.8
$ws.Range("I100").Value = 2049.75
$ws.Range("K100").Value = 2049.75
$ws.Range("M100").Value = -1508.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 8630.75
$ws.Range("I136").Value = 6506.6665
$ws.Range("K136").Value = 19519.9995
$ws.Range("M136").Value = -16969.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3531.3635
$ws.Range("I81").Value = 2245.5454
$ws.Range("K81").Value = 4491.0908
$ws.Range("M81").Value = -3430.0908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3531.3635
$ws.Range("I84").Value = 2245.5454
$ws.Range("K84").Value = 22455.454
$ws.Range("M84").Value = -17151.454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 331.0909
$ws.Range("I113").Value = 277.42856
$ws.Range("K113").Value = 832.28568
$ws.Range("M113").Value = 1337.71432
